# Add the "Singe_Parent_HHs" data worksheet (single parent households table)
# as the last sheet in the workbook, and register it in the TOC sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new worksheet and move it to the end of the tab order -
$ws = $wb.Worksheets.Add()
$ageSheet = $wb.Worksheets.Item("Age")
$ws.Move($null, $ageSheet)

# Re-resolve a live reference to the sheet now that it sits at the end.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Singe_Parent_HHs"

# --- 2. Header row (bold, centered - matches the other data tabs) --------
$newSheet.Cells.Item(1, 1).Value = "single_parent_households"
$newSheet.Cells.Item(1, 2).Value = "county"
$newSheet.Cells.Item(1, 3).Value = "percentages"
$newSheet.Range("A1:C1").Font.Bold = $true
$newSheet.Range("A1:C1").HorizontalAlignment = -4108

# --- 3. Data rows ----------------------------------------------------------
$rows = @(
    @("female_hh", "Imperial", 11.26),
    @("female_hh", "Los Angeles", 6.2),
    @("female_hh", "Orange", 4.7),
    @("female_hh", "Riverside", 5.95),
    @("female_hh", "San Bernardino", 7.62),
    @("female_hh", "Ventura", 4.62),
    @("female_hh", "SCAG", 6.04),
    @("male_hh", "Imperial", 1.6),
    @("male_hh", "Los Angeles", 1.53),
    @("male_hh", "Orange", 1.31),
    @("male_hh", "Riverside", 1.59),
    @("male_hh", "San Bernardino", 1.77),
    @("male_hh", "Ventura", 1.04),
    @("male_hh", "SCAG", 1.51)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# --- 4. Register the new sheet on the TOC page ----------------------------
$toc = $wb.Worksheets.Item("TOC")
$toc.Cells.Item(16, 1).Value = "Singe_Parent_HHs"
$toc.Cells.Item(16, 2).Value = "Single Parent Households (%) by Gender of Head of Household and then by County and SCAG Region"
